# Week 17 data logging: add new WR "S.Williams" row to the WR sheet,
# and make the WR sheet the active/selected tab (was previously QB).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# New player row (row 10) with zeroed stat columns B:J.
$ws.Cells.Item(10, 1).Value = "S.Williams"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0

# Make WR the active sheet/tab (tabSelected moves from QB to WR,
# bookViews.activeTab becomes 2), with the selection left on J11.
$ws.Activate()
$ws.Range("J11").Select()
